# Applies the "Added residual blocks to improve accuracy" edit:
#  - renames sheet "2i" -> "convol_only" and tweaks its view/row height
#  - inserts a new "residual_conn" sheet (right after convol_only) with the
#    num_resid results table, new shared strings and a new red-font style

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) convol_only (formerly "2i")
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "convol_only"

$ws1.Rows.Item(1).RowHeight = 30

# Record the new selection for this sheet *before* a different sheet becomes
# active, so its stored cursor moves from D9 to A3 without leaving this sheet
# marked as the active tab.
$ws1.Range("A3").Select()

# ---------------------------------------------------------------------------
# 2) residual_conn (brand new sheet, placed right after convol_only)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "residual_conn"

# Column widths (matches convol_only's A/B columns; new widths for C/D/E)
$ws2.Columns.Item(1).ColumnWidth = 14.7109375
$ws2.Columns.Item(2).ColumnWidth = 13.5703125
$ws2.Columns.Item(3).ColumnWidth = 14
$ws2.Columns.Item(4).ColumnWidth = 14.85546875
$ws2.Columns.Item(5).ColumnWidth = 52.5703125

# --- seed brand-new shared strings in the exact order they were first used ---
$ws2.Range("A5").Value = "Bayes"
$ws2.Range("D1").Value = "val_accuracy"
$ws2.Range("E7").Value = "Training did not finish due to usage limits"
$ws2.Range("A8").Value = "Manual. 2nd round"
$ws2.Range("C1").Value = "num_resid"
$ws2.Range("E9").Value = "Second attempt"
$ws2.Range("E12").Value = "Second attempt"

# --- remaining header / section-label text (reuses existing shared strings) ---
$ws2.Range("A1").Value = "init_filters"
$ws2.Range("B1").Value = "kernel_size"
$ws2.Range("A2").Value = "Manual"

# --- data rows ---
$ws2.Range("A3").Value = 32
$ws2.Range("B3").Value = 2
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = 0.96153146028518599

$ws2.Range("A4").Value = 32
$ws2.Range("B4").Value = 2
$ws2.Range("C4").Value = 5
$ws2.Range("D4").Value = 0.963708937168121

$ws2.Range("A6").Value = 64
$ws2.Range("B6").Value = 3
$ws2.Range("C6").Value = 5
$ws2.Range("D6").Value = 0.96588641405105502

$ws2.Range("A7").Value = 16
$ws2.Range("B7").Value = 3
$ws2.Range("C7").Value = 6
$ws2.Range("D7").Value = 0.96550000000000002

$ws2.Range("A9").Value = 16
$ws2.Range("B9").Value = 3
$ws2.Range("C9").Value = 6
$ws2.Range("D9").Value = 0.95772093534469604

$ws2.Range("A10").Value = 64
$ws2.Range("B10").Value = 3
$ws2.Range("C10").Value = 7
$ws2.Range("D10").Value = 0.963708937168121

$ws2.Range("A11").Value = 64
$ws2.Range("B11").Value = 4
$ws2.Range("C11").Value = 6
$ws2.Range("D11").Value = 0.96606785058975198

$ws2.Range("A12").Value = 64
$ws2.Range("B12").Value = 4
$ws2.Range("C12").Value = 6
$ws2.Range("D12").Value = 0.96516060829162598

$ws2.Range("A13").Value = 16
$ws2.Range("B13").Value = 5
$ws2.Range("C13").Value = 6
$ws2.Range("D13").Value = 0.96570497751235895

$ws2.Range("A14").Value = 64
$ws2.Range("B14").Value = 5
$ws2.Range("C14").Value = 6
$ws2.Range("D14").Value = 0.96407186985015803

# --- formatting ---

# Header row 1: A1:C1 centered+wrap, D1 centered+wrap+numberformat
$hdr = $ws2.Range("A1:C1")
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4108
$hdr.WrapText = $true

$hdrD = $ws2.Range("D1")
$hdrD.HorizontalAlignment = -4108
$hdrD.VerticalAlignment = -4108
$hdrD.WrapText = $true
$hdrD.NumberFormat = "0.000000"

# Section header rows (2, 5, 8): A:C bold+centered+wrap, D bold+centered+wrap+numberformat
foreach ($r in 2, 5, 8) {
    $sec = $ws2.Range("A" + $r + ":C" + $r)
    $sec.HorizontalAlignment = -4108
    $sec.VerticalAlignment = -4108
    $sec.WrapText = $true
    $sec.Font.Bold = $true

    $secD = $ws2.Range("D" + $r)
    $secD.HorizontalAlignment = -4108
    $secD.VerticalAlignment = -4108
    $secD.WrapText = $true
    $secD.Font.Bold = $true
    $secD.NumberFormat = "0.000000"
}
$ws2.Rows.Item(8).RowHeight = 30

# Plain data cells in columns A:C (all data rows) -> wrap text only
$ws2.Range("A3:C4").WrapText = $true
$ws2.Range("A6:C7").WrapText = $true
$ws2.Range("A9:C14").WrapText = $true

# Column D data values (non-highlighted) -> wrap + numberformat
$dData = $ws2.Range("D3:D4,D6:D7,D9:D10,D13:D14")
$dData.WrapText = $true
$dData.NumberFormat = "0.000000"

# Column D highlighted (red) values -> wrap + numberformat + red font
$dRed = $ws2.Range("D11:D12")
$dRed.WrapText = $true
$dRed.NumberFormat = "0.000000"
$dRed.Font.Color = 255

# Note cells in column E -> wrap text only
$ws2.Range("E7").WrapText = $true
$ws2.Range("E9").WrapText = $true
$ws2.Range("E12").WrapText = $true

# --- page setup (mirrors convol_only) ---
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# --- selection / active sheet ---
$ws2.Range("D12").Select()
